$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 2643.5557
$ws.Range("I4").Value = 2643.5557
$ws.Range("K4").Value = 2643.5557
$ws.Range("M4").Value = -2529.5557
$ws.Range("H18").Value = 3900
$ws.Range("I18").Value = 3350
$ws.Range("K18").Value = 3350
$ws.Range("M18").Value = -3066
$ws.Range("H33").Value = 144.6842
$ws.Range("I33").Value = 148.55556
$ws.Range("K33").Value = 148.55556
$ws.Range("M33").Value = 80.44443999999999
$ws.Range("H38").Value = 522.4286
$ws.Range("I38").Value = 522.4286
$ws.Range("K38").Value = 1567.2858
$ws.Range("M38").Value = -1195.2858
$ws.Range("H43").Value = 5000
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 5000
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 5000
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -5138
$ws.Range("H54").Value = 15500
$ws.Range("J54").Value = 21000
$ws.Range("L54").Value = 21000
$ws.Range("N54").Value = -21972
$ws.Range("H61").Value = 7508.5
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 7508.5
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 22525.5
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -22869.5
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H116").Value = 4999.5
$ws.Range("J116").Value = 4999
$ws.Range("L116").Value = 4999
$ws.Range("N116").Value = -11883
$ws.Range("H127").Value = 613
$ws.Range("I127").Value = 613
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 1839
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = 3121
$ws.Range("N127").ClearContents()
$ws.Range("H138").Value = 13937.25
$ws.Range("I138").Value = 13937.25
$ws.Range("K138").Value = 41811.75
$ws.Range("M138").Value = -36671.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 827.75
$ws.Range("I2").Value = 827.75
$ws.Range("K2").Value = 827.75
$ws.Range("M2").Value = -714.75
$ws.Range("H37").Value = 18000
$ws.Range("I37").Value = 18000
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 18000
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -17727
$ws.Range("N37").ClearContents()
$ws.Range("H46").Value = 3000
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 3000
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 3000
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -3638
$ws.Range("H74").Value = 2419
$ws.Range("I74").Value = 2539.8
$ws.Range("K74").Value = 2539.8
$ws.Range("M74").Value = -1665.8
$ws.Range("H77").Value = 2419
$ws.Range("I77").Value = 2539.8
$ws.Range("K77").Value = 12699
$ws.Range("M77").Value = -8331
$ws.Range("H97").Value = 1719
$ws.Range("I97").Value = 606.6
$ws.Range("K97").Value = 606.6
$ws.Range("M97").Value = -110.6
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H116").Value = 827.75
$ws.Range("I116").Value = 827.75
$ws.Range("K116").Value = 827.75
$ws.Range("M116").Value = 1466.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 827.75
$ws.Range("I3").Value = 827.75
$ws.Range("K3").Value = 827.75
$ws.Range("M3").Value = -713.75
$ws.Range("H20").Value = 3029
$ws.Range("I20").Value = 3029
$ws.Range("K20").Value = 3029
$ws.Range("M20").Value = -2782
$ws.Range("H35").Value = 11333.333
$ws.Range("I35").Value = 4000
$ws.Range("J35").Value = 15000
$ws.Range("K35").Value = 4000
$ws.Range("L35").Value = 15000
$ws.Range("M35").Value = -3690
$ws.Range("N35").Value = -15620
$ws.Range("H100").Value = 16285.25
$ws.Range("J100").Value = 16285.25
$ws.Range("L100").Value = 16285.25
$ws.Range("N100").Value = -18449.25
$ws.Range("H105").Value = 1779.6
$ws.Range("I105").Value = 1474.5
$ws.Range("K105").Value = 1474.5
$ws.Range("M105").Value = 272.5
$ws.Range("H130").Value = 64285.715
$ws.Range("J130").Value = 64285.715
$ws.Range("L130").Value = 64285.715
$ws.Range("N130").Value = -74325.715
$ws.Range("H134").Value = 6935.7646
$ws.Range("J134").Value = 7922.846
$ws.Range("L134").Value = 23768.538
$ws.Range("N134").Value = -28838.538

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2618.889
$ws.Range("I105").Value = 510
$ws.Range("K105").Value = 510
$ws.Range("M105").Value = 1237
$ws.Range("H134").Value = 933.3333
$ws.Range("I134").Value = 933.3333
$ws.Range("K134").Value = 2799.9999
$ws.Range("M134").Value = -264.9998999999998

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 2086.25
$ws.Range("I23").Value = 1939
$ws.Range("J23").Value = 2331.6667
$ws.Range("K23").Value = 5817
$ws.Range("L23").Value = 6995.000100000001
$ws.Range("M23").Value = -5582
$ws.Range("N23").Value = -7465.000100000001
$ws.Range("H34").Value = 1405.1428
$ws.Range("I34").Value = 53.333332
$ws.Range("J34").Value = 2419
$ws.Range("K34").Value = 159.999996
$ws.Range("L34").Value = 7257
$ws.Range("M34").Value = -75.99999600000001
$ws.Range("N34").Value = -7425
$ws.Range("H131").Value = 2403.2942
$ws.Range("I131").Value = 2706.4546
$ws.Range("J131").Value = 1847.5
$ws.Range("K131").Value = 8119.3638
$ws.Range("L131").Value = 5542.5
$ws.Range("M131").Value = -3079.3638
$ws.Range("N131").Value = -15622.5
$ws.Range("H139").Value = 1200
$ws.Range("I139").Value = 1200
$ws.Range("K139").Value = 3600
$ws.Range("M139").Value = 1540

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2900
$ws.Range("J80").Value = 2900
$ws.Range("L80").Value = 2900
$ws.Range("N80").Value = -4896
$ws.Range("H83").Value = 2900
$ws.Range("J83").Value = 2900
$ws.Range("L83").Value = 14500
$ws.Range("N83").Value = -24484
$ws.Range("H97").Value = 1456.4445
$ws.Range("I97").Value = 1058
$ws.Range("K97").Value = 1058
$ws.Range("M97").Value = -562
$ws.Range("H101").Value = 22499.5
$ws.Range("J101").Value = 22499.5
$ws.Range("L101").Value = 22499.5
$ws.Range("N101").Value = -28989.5
$ws.Range("H102").Value = 4779.1665
$ws.Range("I102").Value = 4779.1665
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 4779.1665
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -3157.1665
$ws.Range("N102").ClearContents()
$ws.Range("H126").Value = 4094.4
$ws.Range("I126").Value = 4118
$ws.Range("K126").Value = 12354
$ws.Range("M126").Value = -9884

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2928.1428
$ws.Range("I7").Value = 2253.5454
$ws.Range("J7").Value = 5401.6665
$ws.Range("K7").Value = 2253.5454
$ws.Range("L7").Value = 5401.6665
$ws.Range("M7").Value = -2141.5454
$ws.Range("N7").Value = -5625.6665
$ws.Range("H46").Value = 879.8
$ws.Range("I46").Value = 799.6667
$ws.Range("J46").Value = 1000
$ws.Range("K46").Value = 799.6667
$ws.Range("L46").Value = 1000
$ws.Range("M46").Value = -611.6667
$ws.Range("N46").Value = -1376
$ws.Range("H68").Value = 2883
$ws.Range("J68").Value = 3247.5
$ws.Range("L68").Value = 3247.5
$ws.Range("N68").Value = -4745.5
$ws.Range("H71").Value = 2883
$ws.Range("J71").Value = 3247.5
$ws.Range("L71").Value = 16237.5
$ws.Range("N71").Value = -23725.5
$ws.Range("H103").Value = 30000
$ws.Range("J103").Value = 30000
$ws.Range("L103").Value = 30000
$ws.Range("N103").Value = -32344
$ws.Range("H126").Value = 2928.1428
$ws.Range("I126").Value = 2253.5454
$ws.Range("J126").Value = 5401.6665
$ws.Range("K126").Value = 6760.6362
$ws.Range("L126").Value = 16204.9995
$ws.Range("M126").Value = -4290.6362
$ws.Range("N126").Value = -21144.9995
$ws.Range("H132").Value = 4833.3335
$ws.Range("I132").Value = 5500
$ws.Range("J132").Value = 1500
$ws.Range("K132").Value = 16500
$ws.Range("L132").Value = 4500
$ws.Range("M132").Value = -13970
$ws.Range("N132").Value = -9560

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 55000
$ws.Range("I2").Value = 55000
$ws.Range("K2").Value = 55000
$ws.Range("M2").Value = -54888
$ws.Range("H4").Value = 593.75
$ws.Range("I4").Value = 107.14286
$ws.Range("K4").Value = 107.14286
$ws.Range("M4").Value = 5.857140000000001
$ws.Range("H112").Value = 63834.25
$ws.Range("J112").Value = 63834.25
$ws.Range("L112").Value = 63834.25
$ws.Range("N112").Value = -66788.25
$ws.Range("H122").Value = 335092.66
$ws.Range("I122").Value = 401711.2
$ws.Range("K122").Value = 1205133.6
$ws.Range("M122").Value = -1202683.6
$ws.Range("H126").Value = 759.6667
$ws.Range("I126").Value = 821.6
$ws.Range("J126").Value = 450
$ws.Range("K126").Value = 2464.8
$ws.Range("L126").Value = 1350
$ws.Range("M126").Value = 5.199999999999818
$ws.Range("N126").Value = -6290
$ws.Range("H132").Value = 7255.048
$ws.Range("I132").Value = 5276.077
$ws.Range("J132").Value = 10470.875
$ws.Range("K132").Value = 15828.231
$ws.Range("L132").Value = 31412.625
$ws.Range("M132").Value = -13298.231
$ws.Range("N132").Value = -36472.625
$ws.Range("H133").Value = 82500
$ws.Range("J133").Value = 82500
$ws.Range("L133").Value = 82500
$ws.Range("N133").Value = -92620
